# Round the numeric data in range B2:H13 to the nearest integer.
# (The "Ontpl_" and "Pot_" result files only need to be written to disk
# as integer data.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B2:H13")

foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = [Math]::Round([double]$val, 0)
    }
}
